# Updates the cryptos price-tracker worksheet (Sheet1) with the latest
# scraped values: price (column D), 1h volume/change (column E), and for
# rows 30-31 the coin name/link also changed because the underlying
# ranking reordered (Bittensor now outranks Binance-Peg BSC-USD).
#
# Cells whose new text would otherwise be auto-parsed by Excel as a
# number (and so lose formatting such as trailing zeros, e.g. "1.00",
# "0.1000", or the "350.80" / "0.0000109" style) are entered with a
# leading apostrophe, which is how Excel keeps a numeric-looking entry
# as literal text (quote-prefixed), matching the workbook's existing
# plain-text price column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.589.30"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "2.667.65"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'600.57"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").Value = "'156.58"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.605"
$ws.Range("E8").Value = "  +2.65%  "
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("D10").Value = "'5.92"
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("D11").Value = "'0.398"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "'29.38"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "3.146.35"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "65.380.61"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "2.673.25"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("D20").Value = "'7.56"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "'350.80"
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'69.72"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "'0.0000109"
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("D25").Value = "'9.68"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  -3.73%  "
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("D29").Value = "'8.08"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "'540.61"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("E33").Value = "  -3.76%  "
$ws.Range("D34").Value = "'6.52"
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("E36").Value = "  -2.99%  "
$ws.Range("D37").Value = "'20.39"
$ws.Range("E37").Value = "  -1.65%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'158.99"
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("E40").Value = "  -3.90%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'42.66"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "'165.57"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("E44").Value = "  -2.58%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("E46").Value = "  -5.70%  "
$ws.Range("D47").Value = "'23.03"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'0.646"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("D50").Value = "'0.1000"
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("D51").Value = "'19.92"
$ws.Range("E51").Value = "  +0.46%  "
